$d = $word.ActiveDocument

# Locate the paragraph "Administrator select "add course" menu." — the two
# new steps ("Input user-email." and "Input password.") need to be inserted
# immediately before it.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "*add course*menu*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'add course menu' paragraph"
}

# Insert "Input user-email." as a new list paragraph right before the target.
$target = $d.Paragraphs.Item($targetIndex)
$target.Range.InsertParagraphBefore()
$p1 = $d.Paragraphs.Item($targetIndex)
$p1.Range.Text = "Input user-email."

# Insert "Input password." as a new list paragraph right before the target
# (which is now one slot further down). A trailing placeholder character
# ("X") is appended temporarily so the zero-length bookmark range below can
# be anchored one character before the paragraph mark — placing a
# zero-length range exactly at a paragraph mark is unreliable, so the
# placeholder is removed again right after the bookmark is created.
$targetIndex = $targetIndex + 1
$target2 = $d.Paragraphs.Item($targetIndex)
$target2.Range.InsertParagraphBefore()
$p2 = $d.Paragraphs.Item($targetIndex)
$p2.Range.Text = "Input password.X"

# Move the "_GoBack" bookmark (previously sitting on the last paragraph of
# the document) onto this new "Input password." paragraph, right after its
# text. Re-adding a bookmark with the same name relocates it, so the old
# occurrence near the end of the document disappears automatically.
$p2 = $d.Paragraphs.Item($targetIndex)
$bmPos = $p2.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the placeholder character now that the bookmark is anchored.
$placeholder = $d.Range($bmPos, $bmPos + 1)
$placeholder.Delete()
